$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the empty paragraph right after "...周日" currently carries a
# paragraph-mark run property (<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>)
# inside its <w:pPr>. The target just wants a fully bare <w:p/>.
#
# We locate that paragraph by its (unique, short) text - it is the first
# empty paragraph that immediately follows the paragraph containing "周日".
# ---------------------------------------------------------------------------
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r" -and $i -gt 1) {
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.Text -like "*周日*") {
            $targetIndex = $i
            break
        }
    }
}

if ($targetIndex -ne -1) {
    $p = $d.Paragraphs.Item($targetIndex)
    # Range covering exactly this paragraph's mark (the whole paragraph,
    # since it has no other content).
    $markRng = $d.Range($p.Range.Start, $p.Range.End)
    $blankParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $markRng.InsertXML($blankParaXml)
}

# ---------------------------------------------------------------------------
# Change 2: the paragraph "晴天 今天非常好" gets three more runs appended
# (" " with eastAsia hint, " " plain, "今天天气不错" with eastAsia hint) and
# loses the paragraph-mark run property that used to sit in its <w:pPr>.
# ---------------------------------------------------------------------------
$targetIndex2 = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*晴天 今天非常好*") {
        $targetIndex2 = $i
        break
    }
}

if ($targetIndex2 -ne -1) {
    $p2 = $d.Paragraphs.Item($targetIndex2)
    # Collapsed range right before this paragraph's own mark - i.e. right
    # after the existing text, still "inside" the paragraph.
    $insertAt = $p2.Range.End - 1
    $insRng = $d.Range($insertAt, $insertAt)
    $newRunsXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>今天天气不错</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insRng.InsertXML($newRunsXml)

    # InsertXML above inserted its payload as a *new* paragraph right after
    # the original one (since the payload itself is a full <w:p>...</w:p>).
    # Merge the two paragraphs back into one by deleting the paragraph mark
    # that now separates them - this also drops the old paragraph-mark
    # run-properties (the eastAsia rFonts hint) that lived in the original
    # paragraph's <w:pPr>.
    $d.Range($insertAt, $insertAt + 1).Delete()
}
